$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = 138456
$ws.Range("C8").Value = 758
$ws.Range("D8").Value = 81800
$ws.Range("E8").Value = 52463
$ws.Range("F8").Value = 4288
$ws.Range("G8").Value = 141
$ws.Range("H8").Value = 4193

# Row 16
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 30670
$ws.Range("C16").Value = 564
$ws.Range("D16").Value = 9729
$ws.Range("E16").Value = 19746
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1195

# Row 17
$ws.Range("A17").Value = "Paises Bajos"
$ws.Range("B17").Value = 30449
$ws.Range("C17").Value = 1235
$ws.Range("D17").Value = 250
$ws.Range("E17").Value = 26740
$ws.Range("F17").Value = 1279
$ws.Range("G17").Value = 144
$ws.Range("H17").Value = 3459

# Row 20
$ws.Range("B20").Value = 14568
$ws.Range("C20").Value = 92
$ws.Range("D20").Value = 9704
$ws.Range("E20").Value = 4454
$ws.Range("F20").Value = 227
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 410

# Row 35
$ws.Range("E35").Value = 6744
$ws.Range("F35").Value = 63
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 161

# Row 59
$ws.Range("B59").Value = 2224
$ws.Range("C59").Value = 17
$ws.Range("D59").Value = 269
$ws.Range("E59").Value = 1847
$ws.Range("F59").Value = 71
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 108

# Row 78
$ws.Range("E78").Value = 887
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 6

# Row 82
$ws.Range("A82").Value = "Cuba"
$ws.Range("B82").Value = 923
$ws.Range("C82").Value = 61
$ws.Range("D82").Value = 171
$ws.Range("E82").Value = 721
$ws.Range("F82").Value = 16
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 31

# Row 83
$ws.Range("A83").Value = "Afganistan"
$ws.Range("B83").Value = 906
$ws.Range("C83").Value = 66
$ws.Range("D83").Value = 99
$ws.Range("E83").Value = 777
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 30

# Row 86
$ws.Range("B86").Value = 750
$ws.Range("C86").Value = 15
$ws.Range("D86").Value = 77
$ws.Range("E86").Value = 661
$ws.Range("F86").Value = 8

# Row 114
$ws.Range("B114").Value = 291
$ws.Range("C114").Value = 7
$ws.Range("D114").Value = 169
$ws.Range("E114").Value = 118
$ws.Range("F114").Value = 11

# Row 124
$ws.Range("D124").Value = 38
$ws.Range("E124").Value = 132

# Row 142
$ws.Range("B142").Value = 88
$ws.Range("C142").Value = 3
$ws.Range("D142").Value = 5

# Row 184
$ws.Range("B184").Value = 17
$ws.Range("C184").Value = 1
$ws.Range("E184").Value = 15
